# Add "LastName_NS" and "Address_AS" columns to Sheet1, between the existing
# "Name_AS" (B) and "Age_NI" (C) columns — shifting Age_NI/DOB_ND from C/D to E/F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at C:D (Age_NI/DOB_ND + their data shift right to E:F).
[void]$ws.Columns("C:D").Insert()

# --- New "LastName_NS" column (C) ---
$ws.Range("C1").Value = "LastName_NS"
$ws.Range("C2").Value = "Klark"
$ws.Range("C3").Value = "Kumar"
$ws.Range("C6").Value = "Khanna"
$ws.Range("C5").Value = "Kapoor"
$ws.Range("C4").Value = "Kunte"

# --- New "Address_AS" column (D) ---
$ws.Range("D1").Value = "Address_AS"
$ws.Range("D5").Value = "Bangalore India"
$ws.Range("D2").Value = "Sydney Australia"
$ws.Range("D3").Value = "PutraJaya Malaysia"
$ws.Range("D4").Value = "Maharashtra India"
$ws.Range("D6").Value = "Delhi India"

# Column widths matching the workbook author's manual resize of the new columns.
$ws.Columns("C").ColumnWidth = 13.125
$ws.Columns("D").ColumnWidth = 16.916667

# Restore the active-cell selection to the last edited cell.
[void]$ws.Range("D6").Select()
